$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: move "de " so it precedes {{MESES}} instead of {{AÑO}},
# splitting the original single run into three runs (matching the
# diff's resulting OOXML structure).
# ------------------------------------------------------------------
$oldPhrase = "y su respectivo índice I(%), durante el(los) mes(es) {{MESES}} de {{AÑO}}, de acuerdo con la información preliminar disponible en nuestra base de datos de la estación meteorológica"
$newPhrase = "y su respectivo índice I(%), durante el(los) mes(es) de {{MESES}} {{AÑO}}, de acuerdo con la información preliminar disponible en nuestra base de datos de la estación meteorológica"

$findRange = $d.Content
$findRange.Find.Execute($oldPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $findRange.Start

# Replace the whole phrase's text first (still a single run at this point).
$findRange.Text = $newPhrase

$part1Len = "y su respectivo índice I(%), durante el(los) mes(es) ".Length
$part2Len = "de ".Length

# Toggle Bold on/off over the first chunk to force Word to split a new
# run at that boundary without altering the visible formatting.
$r1 = $d.Range($start, $start + $part1Len)
$r1.Font.Bold = 1
$r1.Font.Bold = 0

# Toggle Bold on/off over the second chunk ("de ") to force another
# run split at its trailing boundary.
$r2 = $d.Range($start + $part1Len, $start + $part1Len + $part2Len)
$r2.Font.Bold = 1
$r2.Font.Bold = 0

# ------------------------------------------------------------------
# Change 2: update the cached TIME field result text.
# ------------------------------------------------------------------
$d.Content.Find.Execute("1 de noviembre de 2024", $true, $false, $false, $false, $false, $true, 1, $false, "17 de diciembre de 2024", 2)
